$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add 6 new "generic" gen-ed-credit course rows (118-123) below the existing
# course table (which currently ends at row 117).
# ---------------------------------------------------------------------------

# Replicate the formatting of the last data row (117) down into the new rows
# for columns A,B,C,E,F,G,H,I,J,K (column D is intentionally left with
# default/no explicit style, matching the source workbook).
$ws.Range("A117:C117").Copy($ws.Range("A118:C123"))
$ws.Range("E117:K117").Copy($ws.Range("E118:K123"))

# id (col A)
$ws.Range("A118").Value = 100009
$ws.Range("A119").Value = 100010
$ws.Range("A120").Value = 100011
$ws.Range("A121").Value = 100012
$ws.Range("A122").Value = 100013
$ws.Range("A123").Value = 100014

# subject (col B) / name (col H) -- written in this interleaved order so the
# resulting shared-string table matches the authored workbook.
$ws.Range("B122").Value = "AAL"
$ws.Range("H122").Value = "Generic Arts and Letters Credit"

$ws.Range("B119").Value = "GP"
$ws.Range("B120").Value = "SCI"
$ws.Range("B121").Value = "SO"
$ws.Range("B118").Value = "US"

$ws.Range("H118").Value = "Generic US differences and inequalities Credit"
$ws.Range("H119").Value = "Generic Global Perspectives Credit"
$ws.Range("H120").Value = "Generic Science Credit"
$ws.Range("H121").Value = "Generic Social Science Credit"

$ws.Range("B123").Value = "CRE"
$ws.Range("H123").Value = "Generic Credit (any)"

# Remaining columns, identical across all 6 new rows:
#   C number, D prereq, E credits, F term, G annual, I required, J keywords, K satisfyarea
foreach ($r in 118..123) {
  $ws.Cells.Item($r, 3).Value = 0
  $ws.Cells.Item($r, 4).Value = $false
  $ws.Cells.Item($r, 5).Value = 4
  $ws.Cells.Item($r, 6).Value = "FWS"
  $ws.Cells.Item($r, 7).Value = $true
  $ws.Cells.Item($r, 9).Value = $false
  $ws.Cells.Item($r, 10).Value = "~"
  $ws.Cells.Item($r, 11).Value = "na"
  $ws.Rows.Item($r).RowHeight = 15.75
}

# ---------------------------------------------------------------------------
# View-state: the author had scrolled down and selected P84 before saving.
# ---------------------------------------------------------------------------
$win = $excel.Windows.Item(1)
$win.ScrollRow = 71
$win.ScrollColumn = 1
$ws.Range("P84").Select()
